$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 186.78572
$ws.Range("I39").Value = 40
$ws.Range("K39").Value = 120
$ws.Range("M39").Value = 176
$ws.Range("H58").Value = 350
$ws.Range("I58").Value = 350
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 1050
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -900
$ws.Range("H87").Value = 71000
$ws.Range("J87").Value = 71250
$ws.Range("L87").Value = 71250
$ws.Range("N87").Value = -73746
$ws.Range("H90").Value = 71000
$ws.Range("J90").Value = 71250
$ws.Range("L90").Value = 213750
$ws.Range("N90").Value = -226230
$ws.Range("H106").Value = 29413616
$ws.Range("I106").Value = 33335104
$ws.Range("K106").Value = 33335104
$ws.Range("M106").Value = -33334473
$ws.Range("H112").Value = 85016.586
$ws.Range("I112").Value = 843
$ws.Range("J112").Value = 113074.445
$ws.Range("K112").Value = 2529
$ws.Range("L112").Value = 339223.335
$ws.Range("M112").Value = -1421
$ws.Range("N112").Value = -341439.335
$ws.Range("H130").Value = 140365
$ws.Range("J130").Value = 140365
$ws.Range("L130").Value = 140365
$ws.Range("N130").Value = -150405
$ws.Range("H137").Value = 1715.75
$ws.Range("J137").Value = 1524
$ws.Range("L137").Value = 4572
$ws.Range("N137").Value = -9672
$ws.Range("H138").Value = 3265.1904
$ws.Range("I138").Value = 2069.4285
$ws.Range("K138").Value = 6208.2855
$ws.Range("M138").Value = -1068.2855
$ws.Range("H140").Value = 102154.2
$ws.Range("J140").Value = 102692.75
$ws.Range("L140").Value = 102692.75
$ws.Range("N140").Value = -113052.75
$ws.Range("N58").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4307.592
$ws.Range("I32").Value = 3940.9348
$ws.Range("K32").Value = 3940.9348
$ws.Range("M32").Value = -3653.9348
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("H45").Value = 9928.1875
$ws.Range("I45").Value = 10125.467
$ws.Range("K45").Value = 10125.467
$ws.Range("M45").Value = -9748.467000000001
$ws.Range("H61").Value = 3036.9023
$ws.Range("I61").Value = 2703.054
$ws.Range("J61").Value = 6125
$ws.Range("K61").Value = 2703.054
$ws.Range("L61").Value = 6125
$ws.Range("M61").Value = -2491.054
$ws.Range("N61").Value = -6549
$ws.Range("H74").Value = 6558.7646
$ws.Range("I74").Value = 1836.8572
$ws.Range("J74").Value = 28594.334
$ws.Range("K74").Value = 1836.8572
$ws.Range("L74").Value = 28594.334
$ws.Range("M74").Value = -962.8571999999999
$ws.Range("N74").Value = -30342.334
$ws.Range("H77").Value = 6558.7646
$ws.Range("I77").Value = 1836.8572
$ws.Range("J77").Value = 28594.334
$ws.Range("K77").Value = 9184.286
$ws.Range("L77").Value = 142971.67
$ws.Range("M77").Value = -4816.286
$ws.Range("N77").Value = -151707.67
$ws.Range("H136").Value = 3036.9023
$ws.Range("I136").Value = 2703.054
$ws.Range("J136").Value = 6125
$ws.Range("K136").Value = 8109.162
$ws.Range("L136").Value = 18375
$ws.Range("M136").Value = -5559.162
$ws.Range("N136").Value = -23475
$ws.Range("N43").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 419637.38
$ws.Range("I22").Value = 740.53845
$ws.Range("J22").Value = 530773.3
$ws.Range("K22").Value = 740.53845
$ws.Range("L22").Value = 530773.3
$ws.Range("M22").Value = -567.53845
$ws.Range("N22").Value = -531119.3
$ws.Range("H86").Value = 1289.2812
$ws.Range("I86").Value = 1423.0769
$ws.Range("K86").Value = 1423.0769
$ws.Range("M86").Value = -300.0769
$ws.Range("H89").Value = 1289.2812
$ws.Range("I89").Value = 1423.0769
$ws.Range("K89").Value = 7115.3845
$ws.Range("M89").Value = -1499.3845

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 474
$ws.Range("I22").Value = 473.5
$ws.Range("J22").Value = 474.5
$ws.Range("K22").Value = 473.5
$ws.Range("L22").Value = 474.5
$ws.Range("M22").Value = -123.5
$ws.Range("N22").Value = -1174.5
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("H58").Value = 2367.8364
$ws.Range("I58").Value = 2418.1667
$ws.Range("J58").Value = 2272.4736
$ws.Range("K58").Value = 2418.1667
$ws.Range("L58").Value = 2272.4736
$ws.Range("M58").Value = -2215.1667
$ws.Range("N58").Value = -2678.4736
$ws.Range("H134").Value = 11050.244
$ws.Range("I134").Value = 4939.543
$ws.Range("J134").Value = 46696
$ws.Range("K134").Value = 14818.629
$ws.Range("L134").Value = 140088
$ws.Range("M134").Value = -12283.629
$ws.Range("N134").Value = -145158
$ws.Range("H136").Value = 2367.8364
$ws.Range("I136").Value = 2418.1667
$ws.Range("J136").Value = 2272.4736
$ws.Range("K136").Value = 7254.500100000001
$ws.Range("L136").Value = 6817.4208
$ws.Range("M136").Value = -4704.500100000001
$ws.Range("N136").Value = -11917.4208
$ws.Range("M25").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 683.7143
$ws.Range("I5").Value = 345.4
$ws.Range("J5").Value = 1529.5
$ws.Range("K5").Value = 1036.2
$ws.Range("L5").Value = 4588.5
$ws.Range("M5").Value = -924.1999999999998
$ws.Range("N5").Value = -4812.5
$ws.Range("H131").Value = 28685.71
$ws.Range("I131").Value = 1000000
$ws.Range("J131").Value = 2433.973
$ws.Range("K131").Value = 3000000
$ws.Range("L131").Value = 7301.919
$ws.Range("M131").Value = -2994960
$ws.Range("N131").Value = -17381.919
$ws.Range("H135").Value = 683.7143
$ws.Range("I135").Value = 345.4
$ws.Range("J135").Value = 1529.5
$ws.Range("K135").Value = 3108.6
$ws.Range("L135").Value = 13765.5
$ws.Range("M135").Value = -573.5999999999999
$ws.Range("N135").Value = -18835.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2937.5
$ws.Range("I80").Value = 2937.5
$ws.Range("K80").Value = 2937.5
$ws.Range("M80").Value = -1939.5
$ws.Range("H83").Value = 2937.5
$ws.Range("I83").Value = 2937.5
$ws.Range("K83").Value = 14687.5
$ws.Range("M83").Value = -9695.5
$ws.Range("H120").Value = 29998.834
$ws.Range("J120").Value = 29998.834
$ws.Range("L120").Value = 29998.834
$ws.Range("N120").Value = -39674.834
$ws.Range("H126").Value = 20808.154
$ws.Range("I126").Value = 47002
$ws.Range("K126").Value = 141006
$ws.Range("M126").Value = -138536
$ws.Range("H132").Value = 2911.5
$ws.Range("J132").Value = 3531.8
$ws.Range("L132").Value = 10595.4
$ws.Range("N132").Value = -15655.4
$ws.Range("H141").Value = 61314.625
$ws.Range("J141").Value = 61314.625
$ws.Range("L141").Value = 61314.625
$ws.Range("N141").Value = -71674.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2004.6666
$ws.Range("I22").Value = 1499
$ws.Range("J22").Value = 2029.95
$ws.Range("K22").Value = 1499
$ws.Range("L22").Value = 2029.95
$ws.Range("M22").Value = -1204
$ws.Range("N22").Value = -2619.95
$ws.Range("H27").Value = 2004.6666
$ws.Range("I27").Value = 1499
$ws.Range("J27").Value = 2029.95
$ws.Range("K27").Value = 1499
$ws.Range("L27").Value = 2029.95
$ws.Range("M27").Value = -1392
$ws.Range("N27").Value = -2243.95
$ws.Range("H82").Value = 2074.7
$ws.Range("I82").Value = 1946.9333
$ws.Range("K82").Value = 1946.9333
$ws.Range("M82").Value = -1585.9333
$ws.Range("H85").Value = 2074.7
$ws.Range("I85").Value = 1946.9333
$ws.Range("K85").Value = 1946.9333
$ws.Range("M85").Value = -698.9332999999999
$ws.Range("H93").Value = 2660
$ws.Range("I93").Value = 2575
$ws.Range("K93").Value = 2575
$ws.Range("M93").Value = -1327
$ws.Range("H127").Value = 26607.5
$ws.Range("J127").Value = 26607.5
$ws.Range("L127").Value = 26607.5
$ws.Range("N127").Value = -36527.5
$ws.Range("H132").Value = 5514.2
$ws.Range("I132").Value = 5126.8887
$ws.Range("K132").Value = 15380.6661
$ws.Range("M132").Value = -12850.6661

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 169733.17
$ws.Range("I3").Value = 501000
$ws.Range("J3").Value = 4099.75
$ws.Range("K3").Value = 501000
$ws.Range("L3").Value = 4099.75
$ws.Range("M3").Value = -500886
$ws.Range("N3").Value = -4327.75
$ws.Range("H81").Value = 10547.934
$ws.Range("I81").Value = 26054.75
$ws.Range("K81").Value = 52109.5
$ws.Range("M81").Value = -51048.5
$ws.Range("H84").Value = 10547.934
$ws.Range("I84").Value = 26054.75
$ws.Range("K84").Value = 260547.5
$ws.Range("M84").Value = -255243.5
$ws.Range("H140").Value = 69257.8
$ws.Range("J140").Value = 69257.8
$ws.Range("L140").Value = 69257.8
$ws.Range("N140").Value = -79617.8
